$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'260.44"
$ws.Range("E2").Value = "'1.89%"
$ws.Range("G2").Value = "'13"
$ws.Range("D3").Value = "'27.31"
$ws.Range("E3").Value = "'3.54%"
$ws.Range("G3").Value = "'13"
$ws.Range("D4").Value = "'4.691"
$ws.Range("E4").Value = "'1.06%"
$ws.Range("G4").Value = "'13"
$ws.Range("D5").Value = "'0.06111"
$ws.Range("E5").Value = "'3.13%"
$ws.Range("G5").Value = "'13"
$ws.Range("D6").Value = "'6.658"
$ws.Range("E6").Value = "'0.50%"
$ws.Range("G6").Value = "'13"
$ws.Range("D7").Value = "'0.8525"
$ws.Range("E7").Value = "'0.14%"
$ws.Range("G7").Value = "'13"
$ws.Range("D8").Value = "'0.9204"
$ws.Range("E8").Value = "'1.90%"
$ws.Range("G8").Value = "'13"
$ws.Range("D9").Value = "'0.1401"
$ws.Range("E9").Value = "'1.62%"
$ws.Range("G9").Value = "'13"
$ws.Range("D10").Value = "'0.04746"
$ws.Range("E10").Value = "'14.71%"
$ws.Range("G10").Value = "'13"
$ws.Range("D11").Value = "'0.07090"
$ws.Range("E11").Value = "'1.47%"
$ws.Range("G11").Value = "'13"
$ws.Range("D12").Value = "'0.03073"
$ws.Range("E12").Value = "'1.38%"
$ws.Range("G12").Value = "'13"
$ws.Range("D13").Value = "'0.09057"
$ws.Range("E13").Value = "'-0.28%"
$ws.Range("G13").Value = "'13"
$ws.Range("D14").Value = "'0.001537"
$ws.Range("E14").Value = "'0.76%"
$ws.Range("G14").Value = "'13"
$ws.Range("D15").Value = "'0.0006100"
$ws.Range("E15").Value = "'0.87%"
$ws.Range("G15").Value = "'13"
$ws.Range("D16").Value = "'0.006017"
$ws.Range("E16").Value = "'-2.47%"
$ws.Range("G16").Value = "'13"
$ws.Range("D17").Value = "'3.453"
$ws.Range("E17").Value = "'-0.50%"
$ws.Range("G17").Value = "'13"
$ws.Range("D18").Value = "'3.147"
$ws.Range("E18").Value = "'0.10%"
$ws.Range("G18").Value = "'13"
$ws.Range("E19").Value = "'-0.63%"
$ws.Range("G19").Value = "'13"
$ws.Range("E20").Value = "'2.95%"
$ws.Range("G20").Value = "'13"
$ws.Range("D21").Value = "'0.1305"
$ws.Range("E21").Value = "'1.54%"
$ws.Range("G21").Value = "'13"
$ws.Range("D22").Value = "'4.100"
$ws.Range("E22").Value = "'6.39%"
$ws.Range("G22").Value = "'13"
$ws.Range("D23").Value = "'0.04225"
$ws.Range("E23").Value = "'0.32%"
$ws.Range("G23").Value = "'13"
$ws.Range("E24").Value = "'0.62%"
$ws.Range("G24").Value = "'13"
$ws.Range("E25").Value = "'-18.72%"
$ws.Range("G25").Value = "'13"
$ws.Range("E26").Value = "'0.04%"
$ws.Range("G26").Value = "'13"
$ws.Range("E27").Value = "'3.41%"
$ws.Range("G27").Value = "'13"
$ws.Range("G28").Value = "'13"
$ws.Range("G29").Value = "'13"
$ws.Range("G30").Value = "'13"
$ws.Range("G31").Value = "'13"
$ws.Range("G32").Value = "'13"
$ws.Range("G33").Value = "'13"
$ws.Range("G34").Value = "'13"
$ws.Range("G35").Value = "'13"
$ws.Range("G36").Value = "'13"
$ws.Range("G37").Value = "'13"
$ws.Range("G38").Value = "'13"
$ws.Range("G39").Value = "'13"
$ws.Range("D40").Value = "'0.03855"
$ws.Range("E40").Value = "'2.25%"
$ws.Range("G40").Value = "'13"
$ws.Range("D41").Value = "'0.1113"
$ws.Range("E41").Value = "'1.80%"
$ws.Range("G41").Value = "'13"
$ws.Range("D42").Value = "'0.004096"
$ws.Range("E42").Value = "'-34.12%"
$ws.Range("G42").Value = "'13"
$ws.Range("D43").Value = "'0.01632"
$ws.Range("E43").Value = "'13.22%"
$ws.Range("G43").Value = "'13"
$ws.Range("E44").Value = "'0.80%"
$ws.Range("G44").Value = "'13"
$ws.Range("D45").Value = "'0.00005157"
$ws.Range("E45").Value = "'0.11%"
$ws.Range("G45").Value = "'13"
$ws.Range("E46").Value = "'0.05%"
$ws.Range("G46").Value = "'13"
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").Value = "'0.05400"
$ws.Range("E47").Value = "'35.06%"
$ws.Range("G47").Value = "'13"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.1355"
$ws.Range("E48").Value = "'-43.74%"
$ws.Range("G48").Value = "'13"
$ws.Range("E49").Value = "'0.05%"
$ws.Range("G49").Value = "'13"
$ws.Range("E50").Value = "'0.05%"
$ws.Range("G50").Value = "'13"
$ws.Range("G51").Value = "'13"
